$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the uniform prior bounds for F4 (row 5), F5 (row 6), F6 (row 7)
$ws.Range("G5").Value = 3
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 5
$ws.Range("F7").Value = 5

# Update the active selection to match the saved view state
$ws.Range("C8").Select()
